# Update crypto price list (symbol list refresh as of Fri Dec 16 23:43:24 UTC 2022)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing "text" storage (since the
# sheet stores all data - including numeric-looking prices - as text), and
# then reset the cell style back to Normal so no stray number-format style
# gets attached to the cell (it would otherwise pick up a Text numFmt style).
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value2 = $text
    $range.Style = "Normal"
}

# --- Column D (Price) updates -------------------------------------------------
Set-TextValue $ws.Range("D2")  "230.19"
Set-TextValue $ws.Range("D3")  "22.51"
Set-TextValue $ws.Range("D4")  "5.526"
Set-TextValue $ws.Range("D5")  "0.05554"
Set-TextValue $ws.Range("D7")  "6.495"
Set-TextValue $ws.Range("D8")  "1.148"
Set-TextValue $ws.Range("D9")  "0.7904"
Set-TextValue $ws.Range("D10") "0.1393"
Set-TextValue $ws.Range("D11") "0.07381"
Set-TextValue $ws.Range("D12") "0.03137"
Set-TextValue $ws.Range("D14") "0.09266"
Set-TextValue $ws.Range("D15") "0.001672"
Set-TextValue $ws.Range("D16") "3.263"
Set-TextValue $ws.Range("D17") "0.04729"
Set-TextValue $ws.Range("D18") "0.0005947"
Set-TextValue $ws.Range("D19") "0.006267"
Set-TextValue $ws.Range("D20") "0.005242"
Set-TextValue $ws.Range("D21") "0.001068"
Set-TextValue $ws.Range("D22") "0.0001506"
Set-TextValue $ws.Range("D23") "3.694"
Set-TextValue $ws.Range("D24") "2.192"
Set-TextValue $ws.Range("D26") "0.1295"
Set-TextValue $ws.Range("D27") "0.0006183"
Set-TextValue $ws.Range("D40") "0.04034"
Set-TextValue $ws.Range("D41") "0.007114"
Set-TextValue $ws.Range("D44") "0.008144"
Set-TextValue $ws.Range("D46") "0.00005542"
Set-TextValue $ws.Range("D47") "0.00000000756"
Set-TextValue $ws.Range("D48") "0.6807"
Set-TextValue $ws.Range("D49") "0.09304"
Set-TextValue $ws.Range("D50") "0.00002117"
Set-TextValue $ws.Range("D51") "0.01018"

# --- Rows 42 / 43: CEJI and BKEXToken swapped places (rank reordering) -------
# Row 42 now holds BKEXToken (previously in row 43), with its updated price.
$ws.Range("B42").Value2 = "BKEXToken"
$ws.Range("C42").Value2 = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D42") "0.1034"
$ws.Range("E42").Value2 = "41BKEXTokenBKK"

# Row 43 now holds CEJI (previously in row 42), with its updated price.
$ws.Range("B43").Value2 = "CEJI"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Range("D43") "0.003221"
$ws.Range("E43").Value2 = "42CEJICEJI"
